$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A1").Value = 45436
$ws.Range("D24").Value = 1057
$ws.Range("D25").Value = 1167
$ws.Range("D26").Value = 1341
$ws.Range("D27").Value = 1353
$ws.Range("D28").Value = 1721
$ws.Range("D29").Value = 1986
$ws.Range("D30").Value = 2336
$ws.Range("D31").Value = 2470
$ws.Range("D36").Value = 1440
$ws.Range("D37").Value = 1760
$ws.Range("D38").Value = 1828
$ws.Range("D39").Value = 1790
$ws.Range("D40").Value = 2165
$ws.Range("D41").Value = 2344
$ws.Range("D42").Value = 2735
$ws.Range("D43").Value = 2955
$ws.Range("D44").Value = 3785
$ws.Range("D45").Value = 4370
$ws.Range("D46").Value = 5025
$ws.Range("D47").Value = 5650
$ws.Range("D52").Value = 1535
$ws.Range("D53").Value = 2056
$ws.Range("D54").Value = 2125
$ws.Range("D55").Value = 2300
$ws.Range("D56").Value = 2500
$ws.Range("D57").Value = 2855
$ws.Range("D58").Value = 3140
$ws.Range("D59").Value = 3486
$ws.Range("D60").Value = 4190
$ws.Range("D61").Value = 5150
$ws.Range("D62").Value = 5810
$ws.Range("D63").Value = 6500
$ws.Range("D68").Value = 1940
$ws.Range("D69").Value = 2030
$ws.Range("D70").Value = 2235
$ws.Range("D71").Value = 2570
$ws.Range("D72").Value = 2950
$ws.Range("D73").Value = 3220
$ws.Range("D74").Value = 3631
$ws.Range("D75").Value = 3920
$ws.Range("D76").Value = 5170
$ws.Range("D77").Value = 5940
$ws.Range("D78").Value = 6920
$ws.Range("D79").Value = 7690
$ws.Range("D80").Value = 8274.946
$ws.Range("D85").Value = 2870
$ws.Range("D86").Value = 3340
$ws.Range("D87").Value = 3650
$ws.Range("D88").Value = 4100
$ws.Range("D89").Value = 4810
$ws.Range("D90").Value = 4690
$ws.Range("D91").Value = 6164
$ws.Range("D92").Value = 7044
$ws.Range("D93").Value = 7790
$ws.Range("D94").Value = 8600
$ws.Range("D95").Value = 10450
$ws.Range("D96").Value = 15450.772
$ws.Range("D101").Value = 4210
$ws.Range("D102").Value = 5220
$ws.Range("D103").Value = 5530
$ws.Range("D104").Value = 6380
$ws.Range("D105").Value = 6748
$ws.Range("D106").Value = 8050
$ws.Range("D107").Value = 8700
$ws.Range("D108").Value = 9950
$ws.Range("D109").Value = 11030
$ws.Range("D110").Value = 15914.403
$ws.Range("D115").Value = 5700
$ws.Range("D116").Value = 6945
$ws.Range("D117").Value = 7340
$ws.Range("D118").Value = 8340
$ws.Range("D119").Value = 8700
$ws.Range("D120").Value = 10550
$ws.Range("D121").Value = 12150
$ws.Range("D122").Value = 13540
$ws.Range("D123").Value = 14750
$ws.Range("D124").Value = 21320.762
